$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date update
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value now filled in
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 (old "Contact") becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Delete old duplicate "Contact" row (row 11), shifting rows up
$ws.Rows.Item(11).Delete()

$wb.Save()
